$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert 2 new rows before row 51, shifting old rows 51-56 down to 53-58.
$ws.Rows.Item(51).Resize(2).Insert()

# Fill new row 51 with data (copy constant columns from row 53 which now holds old row51's data)
$ws.Range("A51").Value = 7
$ws.Range("B51").Value = "Terminal Hortofrutícola Agro Chillán"
$ws.Range("C51").Value = "Ñuble"
$ws.Range("D51").Value = 45209
$ws.Range("D51").Style = $ws.Range("D53").Style
$ws.Range("D51").NumberFormat = $ws.Range("D53").NumberFormat
$ws.Range("E51").Value = 16
$ws.Range("F51").Value = 300000000
$ws.Range("G51").Value = "Espárragos"
$ws.Range("H51").Value = "Sin especificar"
$ws.Range("I51").Value = "Primera"
$ws.Range("J51").Value = 400
$ws.Range("K51").Value = 1200
$ws.Range("L51").Value = 1300
$ws.Range("M51").Value = 1250
$ws.Range("N51").Value = "$/kilo"
$ws.Range("O51").Value = "Región de Ñuble"
$ws.Range("P51").Value = 1250
$ws.Range("Q51").Value = 1
$ws.Range("R51").Value = "Hortaliza"

# Fill new row 52 with data
$ws.Range("A52").Value = 7
$ws.Range("B52").Value = "Terminal Hortofrutícola Agro Chillán"
$ws.Range("C52").Value = "Ñuble"
$ws.Range("D52").Value = 45209
$ws.Range("D52").Style = $ws.Range("D53").Style
$ws.Range("D52").NumberFormat = $ws.Range("D53").NumberFormat
$ws.Range("E52").Value = 16
$ws.Range("F52").Value = 300000000
$ws.Range("G52").Value = "Espárragos"
$ws.Range("H52").Value = "Sin especificar"
$ws.Range("I52").Value = "Primera"
$ws.Range("J52").Value = 300
$ws.Range("K52").Value = 1500
$ws.Range("L52").Value = 1500
$ws.Range("M52").Value = 1500
$ws.Range("N52").Value = "$/kilo"
$ws.Range("O52").Value = "Región del Maule"
$ws.Range("P52").Value = 1500
$ws.Range("Q52").Value = 1
$ws.Range("R52").Value = "Hortaliza"
